# Button theme for till map now customises and transfers between runs.
# Rework the "Product Sheet" data rows: clear out the old rows and write
# the new set of rows/values (products, prices, theme button names and
# foreground colours).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Sheet")

# Clear the old data rows entirely (rows 2 through 40 covers all prior
# entries at rows 11, 21, 22, 28, 29, 34) so stale cells left at rows that
# are no longer used (11, 21, 28, 29) don't linger.
$ws.Range("A2:E40").Clear()

# Re-write the data rows with the new product list. The PRICE column
# (B) is stored as text in this workbook (matching the other product
# rows), so numeric-looking entries get a leading apostrophe to keep
# them as text instead of being interpreted as numbers.
$ws.Range("A22").Value = "ORANGE"
$ws.Range("B22").Value = "'33"
$ws.Range("C22").Value = "btnDefaultItemTheme"
$ws.Range("D22").Value = "Black"

$ws.Range("A23").Value = "PIE"
$ws.Range("B23").Value = "'22"
$ws.Range("C23").Value = "btnLimeGreenItemTheme"
$ws.Range("D23").Value = "DarkGreen"

$ws.Range("A27").Value = "GRAVY"
$ws.Range("B27").Value = "'11"
$ws.Range("C27").Value = "btnDefaultItemTheme"
$ws.Range("D27").Value = "Black"

$ws.Range("A33").Value = "POOP"
$ws.Range("B33").Value = "POOP"
$ws.Range("C33").Value = "btnLimeGreenItemTheme"
$ws.Range("D33").Value = "Green"

$ws.Range("A34").Value = "PIE"
$ws.Range("B34").Value = "PIE"
$ws.Range("C34").Value = "btnDefaultItemTheme"
$ws.Range("D34").Value = "Black"
